$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Structural changes
# ------------------------------------------------------------------

# Drop the now-superfluous 11th column (K). The <cols> width/style
# definitions only ever covered columns A-J, so this does not disturb them.
$ws.Columns.Item(11).Delete()

# Insert a new blank row above the header row (row 2) for the new
# "年级" (grade) field.
$ws.Rows.Item(2).Insert()

# Compensate so the sheet keeps the same total row count (28): remove one
# of the now-shifted blank data rows (the former row 3, now row 4).
$ws.Rows.Item(4).Delete()

# ------------------------------------------------------------------
# 2) Cell values
# ------------------------------------------------------------------

# Row 2 (new "年级" row) - label only, rest of the row stays blank.
$ws.Range("A2").Value = "年级"

# Row 3 (column headers, shifted down from the old row 2): 开课年级 is
# dropped, and 人数 / 任课教师 / 备注 shift one column to the left.
$ws.Range("H3").Value = "人数"
$ws.Range("I3").Value = "任课教师"
$ws.Range("J3").Value = "备注"

# ------------------------------------------------------------------
# 3) Styling fix-ups
# ------------------------------------------------------------------

# Row 2 styling: same look as the header row (A2 like a row header cell,
# B2:G2/J2 like text header cells, H2:I2 like numeric header cells).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B2:G2").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)
$ws.Range("J3").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# Row 3 / I3 ("任课教师") should look like the other text headers, not the
# numeric header it inherited from the old I2 ("人数").
$ws.Range("J3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("J3").Value = "备注"

# Row 4 / H4 (blank data cell) should look like a normal text data cell,
# not the numeric one it inherited from the old H3.
$ws.Range("I4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 1 title bar: give the rest of the merged title row a visible top and
# bottom border; the last cell (J1, end of the merged band) also gets a
# right border to close off the box.
$ws.Range("B1:I1").Borders.Item(8).LineStyle = 1
$ws.Range("B1:I1").Borders.Item(9).LineStyle = 1
$ws.Range("J1").Borders.Item(8).LineStyle = 1
$ws.Range("J1").Borders.Item(9).LineStyle = 1
$ws.Range("J1").Borders.Item(10).LineStyle = 1

# ------------------------------------------------------------------
# 4) Selection / misc
# ------------------------------------------------------------------
$ws.Range("J4").Select()
